$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column at G ("单位简称") - shifts old G..J (类型/地区/联系人/联系电话) to H..K.
#    Style/border of the header row is inherited automatically from the neighboring cell.
$ws.Columns("G:G").Insert()

# 2) Insert a new row at position 7 for the new "天津大学..." record, shifting old rows 7-16 down to 8-17.
$ws.Rows("7:7").Insert()

# 3) Column A holds numbers that must stay stored as *text* (e.g. "1", "2" ... "6"),
#    matching the original workbook convention. Mark the format as Text first so the
#    numeric-looking values are not auto-converted to real numbers, then strip the
#    number-format override again so no stray style id is left on the cells.
$ws.Range("A2:A17").NumberFormat = "@"

# --- Row 1 header ---
$ws.Range("G1").Value = "单位简称"

# --- Data rows 2-17 ---
# row 2
$ws.Range("A2").Value = "1"
$ws.Range("B2").Value = "GX-SB-ZJDX"
$ws.Range("C2").Value = "展品1"
$ws.Range("D2").Value = "测试1"
$ws.Range("E2").Value = "实物1"
$ws.Range("F2").Value = "浙江大学"
$ws.Range("G2").Value = "ZJDX"
$ws.Range("H2").Value = "高校院所"
$ws.Range("I2").Value = "省本级"
$ws.Range("J2").Value = "name"
$ws.Range("K2").Value = "tel135"

# row 3
$ws.Range("A3").Value = "2"
$ws.Range("B3").Value = "GX-SB-XHDX"
$ws.Range("C3").Value = "展品2"
$ws.Range("D3").Value = "测试2"
$ws.Range("E3").Value = "实物2"
$ws.Range("F3").Value = "西湖大学"
$ws.Range("G3").Value = "XHDX"
$ws.Range("H3").Value = "高校院所"
$ws.Range("I3").Value = "省本级"
$ws.Range("J3").Value = "name"
$ws.Range("K3").Value = "tel135"

# row 4
$ws.Range("A4").Value = "3"
$ws.Range("B4").Value = "GX-SB-GYDX"
$ws.Range("C4").Value = "展品3"
$ws.Range("D4").Value = "测试3"
$ws.Range("E4").Value = "实物3"
$ws.Range("F4").Value = "浙江工业大学"
$ws.Range("G4").Value = "GYDX"
$ws.Range("H4").Value = "高校院所"
$ws.Range("I4").Value = "省本级"
$ws.Range("J4").Value = "name"
$ws.Range("K4").Value = "tel135"

# row 5
$ws.Range("A5").Value = "4"
$ws.Range("B5").Value = "GX-SB-ZJSFDX"
$ws.Range("C5").Value = "展品4"
$ws.Range("D5").Value = "测试4"
$ws.Range("E5").Value = "实物4"
$ws.Range("F5").Value = "浙江师范大学"
$ws.Range("G5").Value = "ZJSFDX"
$ws.Range("H5").Value = "高校院所"
$ws.Range("I5").Value = "省本级"
$ws.Range("J5").Value = "name"
$ws.Range("K5").Value = "tel135"

# row 6
$ws.Range("A6").Value = "5"
$ws.Range("B6").Value = "GX-SB-NBDX"
$ws.Range("C6").Value = "展品5"
$ws.Range("D6").Value = "测试5"
$ws.Range("E6").Value = "实物5"
$ws.Range("F6").Value = "宁波大学"
$ws.Range("G6").Value = "NBDX"
$ws.Range("H6").Value = "高校院所"
$ws.Range("I6").Value = "省本级"
$ws.Range("J6").Value = "name"
$ws.Range("K6").Value = "tel135"

# row 7
$ws.Range("A7").Value = "6"
$ws.Range("B7").Value = "ZYJY-SX-TDSX"
$ws.Range("F7").Value = "天津大学浙江绍兴研究院"
$ws.Range("G7").Value = "TDSX"
$ws.Range("H7").Value = "省实验室和省级新型研发机构"
$ws.Range("I7").Value = "绍兴市"

# row 8
$ws.Range("A8").Value = "1"
$ws.Range("B8").Value = "GX-SB-ZJDX"
$ws.Range("C8").Value = "展品1"
$ws.Range("D8").Value = "测试1"
$ws.Range("E8").Value = "实物1"
$ws.Range("F8").Value = "浙江大学"
$ws.Range("G8").Value = "ZJDX"
$ws.Range("H8").Value = "高校院所"
$ws.Range("I8").Value = "省本级"
$ws.Range("J8").Value = "name"
$ws.Range("K8").Value = "tel135"

# row 9
$ws.Range("A9").Value = "2"
$ws.Range("B9").Value = "GX-SB-XHDX"
$ws.Range("C9").Value = "展品2"
$ws.Range("D9").Value = "测试2"
$ws.Range("E9").Value = "实物2"
$ws.Range("F9").Value = "西湖大学"
$ws.Range("G9").Value = "XHDX"
$ws.Range("H9").Value = "高校院所"
$ws.Range("I9").Value = "省本级"
$ws.Range("J9").Value = "name"
$ws.Range("K9").Value = "tel135"

# row 10
$ws.Range("A10").Value = "3"
$ws.Range("B10").Value = "GX-SB-GYDX"
$ws.Range("C10").Value = "展品3"
$ws.Range("D10").Value = "测试3"
$ws.Range("E10").Value = "实物3"
$ws.Range("F10").Value = "浙江工业大学"
$ws.Range("G10").Value = "GYDX"
$ws.Range("H10").Value = "高校院所"
$ws.Range("I10").Value = "省本级"
$ws.Range("J10").Value = "name"
$ws.Range("K10").Value = "tel135"

# row 11
$ws.Range("A11").Value = "4"
$ws.Range("B11").Value = "GX-SB-ZJSFDX"
$ws.Range("C11").Value = "展品4"
$ws.Range("D11").Value = "测试4"
$ws.Range("E11").Value = "实物4"
$ws.Range("F11").Value = "浙江师范大学"
$ws.Range("G11").Value = "ZJSFDX"
$ws.Range("H11").Value = "高校院所"
$ws.Range("I11").Value = "省本级"
$ws.Range("J11").Value = "name"
$ws.Range("K11").Value = "tel135"

# row 12
$ws.Range("A12").Value = "5"
$ws.Range("B12").Value = "GX-SB-NBDX"
$ws.Range("C12").Value = "展品5"
$ws.Range("D12").Value = "测试5"
$ws.Range("E12").Value = "实物5"
$ws.Range("F12").Value = "宁波大学"
$ws.Range("G12").Value = "NBDX"
$ws.Range("H12").Value = "高校院所"
$ws.Range("I12").Value = "省本级"
$ws.Range("J12").Value = "name"
$ws.Range("K12").Value = "tel135"

# row 13
$ws.Range("A13").Value = "1"
$ws.Range("B13").Value = "GX-SB-ZJDX"
$ws.Range("C13").Value = "展品1"
$ws.Range("D13").Value = "测试1"
$ws.Range("E13").Value = "实物1"
$ws.Range("F13").Value = "浙江大学"
$ws.Range("G13").Value = "ZJDX"
$ws.Range("H13").Value = "高校院所"
$ws.Range("I13").Value = "省本级"
$ws.Range("J13").Value = "name"
$ws.Range("K13").Value = "tel135"

# row 14
$ws.Range("A14").Value = "2"
$ws.Range("B14").Value = "GX-SB-XHDX"
$ws.Range("C14").Value = "展品2"
$ws.Range("D14").Value = "测试2"
$ws.Range("E14").Value = "实物2"
$ws.Range("F14").Value = "西湖大学"
$ws.Range("G14").Value = "XHDX"
$ws.Range("H14").Value = "高校院所"
$ws.Range("I14").Value = "省本级"
$ws.Range("J14").Value = "name"
$ws.Range("K14").Value = "tel135"

# row 15
$ws.Range("A15").Value = "3"
$ws.Range("B15").Value = "GX-SB-GYDX"
$ws.Range("C15").Value = "展品3"
$ws.Range("D15").Value = "测试3"
$ws.Range("E15").Value = "实物3"
$ws.Range("F15").Value = "浙江工业大学"
$ws.Range("G15").Value = "GYDX"
$ws.Range("H15").Value = "高校院所"
$ws.Range("I15").Value = "省本级"
$ws.Range("J15").Value = "name"
$ws.Range("K15").Value = "tel135"

# row 16
$ws.Range("A16").Value = "4"
$ws.Range("B16").Value = "GX-SB-ZJSFDX"
$ws.Range("C16").Value = "展品4"
$ws.Range("D16").Value = "测试4"
$ws.Range("E16").Value = "实物4"
$ws.Range("F16").Value = "浙江师范大学"
$ws.Range("G16").Value = "ZJSFDX"
$ws.Range("H16").Value = "高校院所"
$ws.Range("I16").Value = "省本级"
$ws.Range("J16").Value = "name"
$ws.Range("K16").Value = "tel135"

# row 17
$ws.Range("A17").Value = "5"
$ws.Range("B17").Value = "GX-SB-NBDX"
$ws.Range("C17").Value = "展品5"
$ws.Range("D17").Value = "测试5"
$ws.Range("E17").Value = "实物5"
$ws.Range("F17").Value = "宁波大学"
$ws.Range("G17").Value = "NBDX"
$ws.Range("H17").Value = "高校院所"
$ws.Range("I17").Value = "省本级"
$ws.Range("J17").Value = "name"
$ws.Range("K17").Value = "tel135"

# 4) Remove the temporary Text number-format now that the values are locked in as text,
#    restoring column A to the default (unstyled) appearance used throughout the sheet.
$ws.Range("A2:A17").ClearFormats()
